$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 44; existing rows 44-82 shift down to 45-83.
$ws.Rows("44:44").Insert()

# Populate the new row 44 (mirrors the "1-4 足球高手" entries pattern).
$ws.Range("A44").Value = "第一冊"
$ws.Range("B44").Value = "CH1"
$ws.Range("C44").Value = "1-4"
$ws.Range("D44").Value = "1-4 足球高手--康軒有GO補"
$ws.Range("E44").Value = "https://xtjh-yucc.github.io/math/071math/1-4/html5_football.html"

# Match row formatting used by sibling rows (custom row height).
$ws.Rows("44:44").RowHeight = 18.6

# Add the hyperlink for the new URL cell.
$ws.Hyperlinks.Add($ws.Range("E44"), "https://xtjh-yucc.github.io/math/071math/1-4/html5_football.html")

# Re-apply the plain (non-wrapping) hyperlink cell style used by sibling
# rows (E46 already carries it) so the hyperlink Add above doesn't leave a
# duplicated style behind.
$ws.Range("E46").Copy()
$ws.Range("E44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the post-edit view state.
$ws.Range("E44").Select()
